$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix invalid facility utilisation ("Facility limit") values in the report fixture
$ws.Range("E2").Value = 800000
$ws.Range("E3").Value = 800000
$ws.Range("E4").Value = 800000
$ws.Range("E5").Value = 800000
$ws.Range("E6").Value = 800000

$ws.Range("G5").Value = 456
$ws.Range("H5").Value = 3938753.8

$ws.Range("G6").Value = 761579.37

# Columns E:H now share the same (best-fit) width
$ws.Range("E1:H1").ColumnWidth = 15.5

# Update the active selection to match the reviewed range
$ws.Range("E2:H6").Select() | Out-Null
